$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Mapping of row number -> (new DAMSLTag, new DialogAct)
# Column I = DAMSLTag, Column J = DialogAct
$changes = @{
    13  = @("sd", "Statement-non-opinion")
    16  = @("sd", "Statement-non-opinion")
    20  = @("sv", "Statement-opinion")
    41  = @("%", "Uninterpretable")
    50  = @("qy", "Yes-No-Question")
    61  = @("sd", "Statement-non-opinion")
    71  = @("aa", "Agree/Accept")
    75  = @("b", "Acknowledge (Backchannel)")
    84  = @("sd", "Statement-non-opinion")
    131 = @("b", "Acknowledge (Backchannel)")
    142 = @("ba", "Appreciation")
    176 = @("b", "Acknowledge (Backchannel)")
    177 = @("b", "Acknowledge (Backchannel)")
    200 = @("b", "Acknowledge (Backchannel)")
    206 = @("sv", "Statement-opinion")
    207 = @("sd", "Statement-non-opinion")
    233 = @("%", "Uninterpretable")
    247 = @("sd", "Statement-non-opinion")
    257 = @("sd", "Statement-non-opinion")
    259 = @("sd", "Statement-non-opinion")
    266 = @("aa", "Agree/Accept")
    275 = @("ba", "Appreciation")
    284 = @("sd", "Statement-non-opinion")
}

foreach ($row in $changes.Keys) {
    $vals = $changes[$row]
    $ws.Range("I$row").Value = $vals[0]
    $ws.Range("J$row").Value = $vals[1]
}

$wb.Save()
